# إضافة حدث جديد في Card17
# Row 13's previously-blank measurement columns (B:K, M) are filled with "nan"
# placeholders, and a brand-new service-log row (14) is appended with a new
# date, correction note and technician name, mirroring the layout of row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card17")

# Duplicate row 13 (in its original, pre-edit state) into row 14 first so the
# new row inherits the same cell typing/formatting as the existing blank
# cells (B:K, M) as well as the text-typed "17" in column A.
$ws.Range("A13:O13").Copy()
$ws.Range("A14:O14").PasteSpecial()

# Fill the previously-empty columns of row 13 with "nan" placeholders.
$ws.Range("B13").Value = "nan"
$ws.Range("C13").Value = "nan"
$ws.Range("D13").Value = "nan"
$ws.Range("E13").Value = "nan"
$ws.Range("F13").Value = "nan"
$ws.Range("G13").Value = "nan"
$ws.Range("H13").Value = "nan"
$ws.Range("I13").Value = "nan"
$ws.Range("J13").Value = "nan"
$ws.Range("K13").Value = "nan"
$ws.Range("M13").Value = "nan"

# Populate the new row 14 with the new service event details.
# (A14 already holds "17", copied from row 13; B14:K14 and M14 stay blank.)
$ws.Range("L14").Value = "9\12\2024"
$ws.Range("N14").Value = "تم سن الفلاتس + صينه نصف سنويه"
$ws.Range("O14").Value = "الخبير"
